$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 44745.7
$ws.Range("I6").Value = 397.66666
$ws.Range("J6").Value = 63752
$ws.Range("K6").Value = 1192.99998
$ws.Range("L6").Value = 191256
$ws.Range("M6").Value = -1080.99998
$ws.Range("N6").Value = -191480

$ws.Range("H116").Value = 6797.778
$ws.Range("I116").Value = 1445
$ws.Range("J116").Value = 11080
$ws.Range("K116").Value = 1445
$ws.Range("L116").Value = 11080
$ws.Range("M116").Value = 1997
$ws.Range("N116").Value = -17964

$ws.Range("H125").Value = 633
$ws.Range("I125").Value = 644
$ws.Range("J125").Value = 600
$ws.Range("K125").Value = 5796
$ws.Range("L125").Value = 5400
$ws.Range("M125").Value = -3336
$ws.Range("N125").Value = -10320

$ws.Range("H138").Value = 3989218.8
$ws.Range("I138").Value = 1309.5
$ws.Range("J138").Value = 5625284
$ws.Range("K138").Value = 3928.5
$ws.Range("L138").Value = 16875852
$ws.Range("M138").Value = 1211.5
$ws.Range("N138").Value = -16886132

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 488.66666
$ws.Range("I80").Value = 119.5
$ws.Range("J80").Value = 545.46155
$ws.Range("K80").Value = 119.5
$ws.Range("L80").Value = 545.46155
$ws.Range("M80").Value = 878.5
$ws.Range("N80").Value = -2541.46155

$ws.Range("H83").Value = 488.66666
$ws.Range("I83").Value = 119.5
$ws.Range("J83").Value = 545.46155
$ws.Range("K83").Value = 597.5
$ws.Range("L83").Value = 2727.30775
$ws.Range("M83").Value = 4394.5
$ws.Range("N83").Value = -12711.30775

$ws.Range("H141").Value = 30613.334
$ws.Range("J141").Value = 30613.334
$ws.Range("L141").Value = 30613.334
$ws.Range("N141").Value = -40973.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 30305108
$ws.Range("I16").Value = 1571.8334
$ws.Range("J16").Value = 66669350
$ws.Range("K16").Value = 1571.8334
$ws.Range("L16").Value = 66669350
$ws.Range("M16").Value = -1284.8334
$ws.Range("N16").Value = -66669924

$ws.Range("H31").Value = 2268.434
$ws.Range("I31").Value = 1123.5128
$ws.Range("J31").Value = 5457.857
$ws.Range("K31").Value = 1123.5128
$ws.Range("L31").Value = 5457.857
$ws.Range("M31").Value = -828.5128
$ws.Range("N31").Value = -6047.857

$ws.Range("H34").Value = 2268.434
$ws.Range("I34").Value = 1123.5128
$ws.Range("J34").Value = 5457.857
$ws.Range("K34").Value = 1123.5128
$ws.Range("L34").Value = 5457.857
$ws.Range("M34").Value = -921.5128
$ws.Range("N34").Value = -5861.857

$ws.Range("H105").Value = 1426
$ws.Range("I105").Value = 810
$ws.Range("K105").Value = 810
$ws.Range("M105").Value = 937

$ws.Range("H113").Value = 30305108
$ws.Range("I113").Value = 1571.8334
$ws.Range("J113").Value = 66669350
$ws.Range("K113").Value = 1571.8334
$ws.Range("L113").Value = 66669350
$ws.Range("M113").Value = 598.1666
$ws.Range("N113").Value = -66673690

$ws.Range("H132").Value = 40916.383
$ws.Range("I132").Value = 3057.4
$ws.Range("J132").Value = 64578.25
$ws.Range("K132").Value = 9172.200000000001
$ws.Range("L132").Value = 193734.75
$ws.Range("M132").Value = -6642.200000000001
$ws.Range("N132").Value = -198794.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3917.75
$ws.Range("I7").Value = 5133.3335
$ws.Range("K7").Value = 15400.0005
$ws.Range("M7").Value = -15288.0005

$ws.Range("H58").Value = 900
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 900
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 2700
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -2956

$ws.Range("H80").Value = 4024.75
$ws.Range("J80").Value = 4024.75
$ws.Range("L80").Value = 12074.25
$ws.Range("N80").Value = -13946.25

$ws.Range("H83").Value = 4024.75
$ws.Range("J83").Value = 4024.75
$ws.Range("L83").Value = 36222.75
$ws.Range("N83").Value = -45582.75

$ws.Range("H92").Value = 1024.25
$ws.Range("I92").Value = 299.66666
$ws.Range("J92").Value = 1459
$ws.Range("K92").Value = 898.9999799999999
$ws.Range("L92").Value = 4377
$ws.Range("M92").Value = 349.0000200000001
$ws.Range("N92").Value = -6873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4000
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -1849
$ws.Range("N43").Value = -8302

$ws.Range("H46").Value = 25046
$ws.Range("J46").Value = 25046
$ws.Range("L46").Value = 25046
$ws.Range("N46").Value = -25358

$ws.Range("H80").Value = 2982.5862
$ws.Range("I80").Value = 2265.3333
$ws.Range("J80").Value = 3751.0715
$ws.Range("K80").Value = 2265.3333
$ws.Range("L80").Value = 3751.0715
$ws.Range("M80").Value = -1267.3333
$ws.Range("N80").Value = -5747.0715

$ws.Range("H83").Value = 2982.5862
$ws.Range("I83").Value = 2265.3333
$ws.Range("J83").Value = 3751.0715
$ws.Range("K83").Value = 11326.6665
$ws.Range("L83").Value = 18755.3575
$ws.Range("M83").Value = -6334.666499999999
$ws.Range("N83").Value = -28739.3575

$ws.Range("H122").Value = 1629.5714
$ws.Range("I122").Value = 1484.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4453.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2003.5
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2092.3
$ws.Range("I40").Value = 2008.2941
$ws.Range("J40").Value = 2568.3333
$ws.Range("K40").Value = 2008.2941
$ws.Range("L40").Value = 2568.3333
$ws.Range("M40").Value = -1872.2941
$ws.Range("N40").Value = -2840.3333

$ws.Range("H61").Value = 5466.6665
$ws.Range("I61").Value = 5466.6665
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5466.6665
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5264.6665
$ws.Range("N61").ClearContents()

$ws.Range("H93").Value = 1196.7059
$ws.Range("I93").Value = 980.8461
$ws.Range("J93").Value = 1898.25
$ws.Range("K93").Value = 980.8461
$ws.Range("L93").Value = 1898.25
$ws.Range("M93").Value = 267.1539
$ws.Range("N93").Value = -4394.25

$ws.Range("H100").Value = 47268.184
$ws.Range("I100").Value = 85241.664
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 85241.664
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -84700.664
$ws.Range("N100").Value = -2782

$ws.Range("H113").Value = 5466.6665
$ws.Range("I113").Value = 5466.6665
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5466.6665
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3296.6665
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 40897.742
$ws.Range("I132").Value = 3447.077
$ws.Range("K132").Value = 10341.231
$ws.Range("M132").Value = -7811.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1469.0714
$ws.Range("I126").Value = 1036.8422
$ws.Range("J126").Value = 2381.5557
$ws.Range("K126").Value = 3110.5266
$ws.Range("L126").Value = 7144.6671
$ws.Range("M126").Value = -640.5266000000001
$ws.Range("N126").Value = -12084.6671

$ws.Range("H132").Value = 102063.2
$ws.Range("I132").Value = 112134.664
$ws.Range("J132").Value = 93822.91
$ws.Range("K132").Value = 336403.992
$ws.Range("L132").Value = 281468.73
$ws.Range("M132").Value = -333873.992
$ws.Range("N132").Value = -286528.73
